$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 12
$ws_ALC.Range("H12").Value = 401
$ws_ALC.Range("I12").Value = 472.25
$ws_ALC.Range("J12").Value = 258.5
$ws_ALC.Range("K12").Value = 472.25
$ws_ALC.Range("L12").Value = 258.5
$ws_ALC.Range("M12").Value = -302.25
$ws_ALC.Range("N12").Value = -598.5

# ALC row 15
$ws_ALC.Range("H15").Value = 723.4677
$ws_ALC.Range("I15").Value = 723.4677
$ws_ALC.Range("K15").Value = 2170.4031
$ws_ALC.Range("M15").Value = -2001.4031

# ALC row 43
$ws_ALC.Range("H43").Value = 5334.8
$ws_ALC.Range("J43").Value = 5969.8
$ws_ALC.Range("L43").Value = 5969.8
$ws_ALC.Range("N43").Value = -6107.8

# ALC row 101
$ws_ALC.Range("H101").Value = 542.9
$ws_ALC.Range("I101").Value = 572.5
$ws_ALC.Range("J101").Value = 424.5
$ws_ALC.Range("K101").Value = 1717.5
$ws_ALC.Range("L101").Value = 1273.5
$ws_ALC.Range("M101").Value = -95.5
$ws_ALC.Range("N101").Value = -4517.5

# ALC row 113
$ws_ALC.Range("H113").Value = 5999.75
$ws_ALC.Range("I113").Value = 12749.5
$ws_ALC.Range("K113").Value = 12749.5
$ws_ALC.Range("M113").Value = -9495.5

# ALC row 132
$ws_ALC.Range("H132").Value = 24400354
$ws_ALC.Range("I132").Value = 27780958
$ws_ALC.Range("K132").Value = 83342874
$ws_ALC.Range("M132").Value = -83340344

# ALC row 135
$ws_ALC.Range("H135").Value = 4580.727
$ws_ALC.Range("I135").Value = 2283.625
$ws_ALC.Range("K135").Value = 20552.625
$ws_ALC.Range("M135").Value = -18017.625

# ALC row 137
$ws_ALC.Range("H137").Value = 3326.4211
$ws_ALC.Range("I137").Value = 3456.5
$ws_ALC.Range("J137").Value = 2632.6667
$ws_ALC.Range("K137").Value = 10369.5
$ws_ALC.Range("L137").Value = 7898.000100000001
$ws_ALC.Range("M137").Value = -7819.5
$ws_ALC.Range("N137").Value = -12998.0001

# ARM row 32
$ws_ARM.Range("H32").Value = 7083.3477
$ws_ARM.Range("I32").Value = 7113.9854
$ws_ARM.Range("K32").Value = 7113.9854
$ws_ARM.Range("M32").Value = -6826.9854

# ARM row 45
$ws_ARM.Range("H45").Value = 2865.96
$ws_ARM.Range("I45").Value = 2554.3333
$ws_ARM.Range("K45").Value = 2554.3333
$ws_ARM.Range("M45").Value = -2177.3333

# ARM row 61
$ws_ARM.Range("H61").Value = 7651.844
$ws_ARM.Range("I61").Value = 6215.48
$ws_ARM.Range("J61").Value = 9447.299999999999
$ws_ARM.Range("K61").Value = 6215.48
$ws_ARM.Range("L61").Value = 9447.299999999999
$ws_ARM.Range("M61").Value = -6003.48
$ws_ARM.Range("N61").Value = -9871.299999999999

# ARM row 63
$ws_ARM.Range("H63").Value = 3150
$ws_ARM.Range("I63").Value = 2300
$ws_ARM.Range("K63").Value = 2300
$ws_ARM.Range("M63").Value = -1614

# ARM row 66
$ws_ARM.Range("H66").Value = 3150
$ws_ARM.Range("I66").Value = 2300
$ws_ARM.Range("K66").Value = 11500
$ws_ARM.Range("M66").Value = -8068

# ARM row 122
$ws_ARM.Range("H122").Value = 3718.628
$ws_ARM.Range("I122").Value = 3331.4243
$ws_ARM.Range("J122").Value = 4996.4
$ws_ARM.Range("K122").Value = 9994.2729
$ws_ARM.Range("L122").Value = 14989.2
$ws_ARM.Range("M122").Value = -7544.2729
$ws_ARM.Range("N122").Value = -19889.2

# ARM row 132
$ws_ARM.Range("H132").Value = 1693.2909
$ws_ARM.Range("I132").Value = 1630.1
$ws_ARM.Range("K132").Value = 4890.299999999999
$ws_ARM.Range("M132").Value = -2360.299999999999

# ARM row 136
$ws_ARM.Range("H136").Value = 7651.844
$ws_ARM.Range("I136").Value = 6215.48
$ws_ARM.Range("J136").Value = 9447.299999999999
$ws_ARM.Range("K136").Value = 18646.44
$ws_ARM.Range("L136").Value = 28341.9
$ws_ARM.Range("M136").Value = -16096.44
$ws_ARM.Range("N136").Value = -33441.89999999999

# BSM row 20
$ws_BSM.Range("H20").Value = 2361.5
$ws_BSM.Range("I20").Value = 2375.4546
$ws_BSM.Range("J20").Value = 2344.4443
$ws_BSM.Range("K20").Value = 2375.4546
$ws_BSM.Range("L20").Value = 2344.4443
$ws_BSM.Range("M20").Value = -2128.4546
$ws_BSM.Range("N20").Value = -2838.4443

# BSM row 105
$ws_BSM.Range("H105").Value = 4648.231
$ws_BSM.Range("I105").Value = 4122.9
$ws_BSM.Range("K105").Value = 4122.9
$ws_BSM.Range("M105").Value = -2375.9

# BSM row 134
$ws_BSM.Range("H134").Value = 5680.3
$ws_BSM.Range("I134").Value = 2322.861
$ws_BSM.Range("K134").Value = 6968.583
$ws_BSM.Range("M134").Value = -4433.583

# CRP row 22
$ws_CRP.Range("H22").Value = 310.52942
$ws_CRP.Range("I22").Value = 301.07144
$ws_CRP.Range("K22").Value = 301.07144
$ws_CRP.Range("M22").Value = 48.92856

# CRP row 31
$ws_CRP.Range("H31").Value = 2086.2554
$ws_CRP.Range("I31").Value = 1542.159
$ws_CRP.Range("K31").Value = 1542.159
$ws_CRP.Range("M31").Value = -1247.159

# CRP row 34
$ws_CRP.Range("H34").Value = 2086.2554
$ws_CRP.Range("I34").Value = 1542.159
$ws_CRP.Range("K34").Value = 1542.159
$ws_CRP.Range("M34").Value = -1340.159

# CRP row 59
$ws_CRP.Range("H59").Value = 43733
$ws_CRP.Range("I59").Value = 25000
$ws_CRP.Range("J59").Value = 53099.5
$ws_CRP.Range("K59").Value = 25000
$ws_CRP.Range("L59").Value = 53099.5
$ws_CRP.Range("M59").Value = -23855
$ws_CRP.Range("N59").Value = -55389.5

# CRP row 105
$ws_CRP.Range("H105").Value = 727
$ws_CRP.Range("I105").Value = 659.125
$ws_CRP.Range("J105").Value = 998.5
$ws_CRP.Range("K105").Value = 659.125
$ws_CRP.Range("L105").Value = 998.5
$ws_CRP.Range("M105").Value = 1087.875
$ws_CRP.Range("N105").Value = -4492.5

# CRP row 132
$ws_CRP.Range("H132").Value = 835032.5600000001
$ws_CRP.Range("I132").Value = 1082429.6
$ws_CRP.Range("J132").Value = 2879
$ws_CRP.Range("K132").Value = 3247288.8
$ws_CRP.Range("L132").Value = 8637
$ws_CRP.Range("M132").Value = -3244758.8
$ws_CRP.Range("N132").Value = -13697

# CRP row 134
$ws_CRP.Range("H134").Value = 4070.05
$ws_CRP.Range("I134").Value = 1679.25
$ws_CRP.Range("J134").Value = 7656.25
$ws_CRP.Range("K134").Value = 5037.75
$ws_CRP.Range("L134").Value = 22968.75
$ws_CRP.Range("M134").Value = -2502.75
$ws_CRP.Range("N134").Value = -28038.75

# CUL row 15
$ws_CUL.Range("H15").Value = 349.5
$ws_CUL.Range("J15").Value = 500
$ws_CUL.Range("L15").Value = 1500
$ws_CUL.Range("N15").Value = -1780

# CUL row 80
$ws_CUL.Range("H80").Value = 2499.5
$ws_CUL.Range("I80").Value = 1000
$ws_CUL.Range("J80").Value = 3999
$ws_CUL.Range("K80").Value = 3000
$ws_CUL.Range("L80").Value = 11997
$ws_CUL.Range("M80").Value = -2064
$ws_CUL.Range("N80").Value = -13869

# CUL row 83
$ws_CUL.Range("H83").Value = 2499.5
$ws_CUL.Range("I83").Value = 1000
$ws_CUL.Range("J83").Value = 3999
$ws_CUL.Range("K83").Value = 9000
$ws_CUL.Range("L83").Value = 35991
$ws_CUL.Range("M83").Value = -4320
$ws_CUL.Range("N83").Value = -45351

# CUL row 107
$ws_CUL.Range("H107").Value = 5680.909
$ws_CUL.Range("I107").Value = 400
$ws_CUL.Range("J107").Value = 6209
$ws_CUL.Range("K107").Value = 1200
$ws_CUL.Range("L107").Value = 18627
$ws_CUL.Range("M107").Value = 720
$ws_CUL.Range("N107").Value = -22467

# CUL row 113
$ws_CUL.Range("H113").Value = 1790.5
$ws_CUL.Range("J113").Value = 1892.3889
$ws_CUL.Range("L113").Value = 5677.1667
$ws_CUL.Range("N113").Value = -10017.1667

# CUL row 122
$ws_CUL.Range("H122").Value = 1078.6875
$ws_CUL.Range("J122").Value = 1108.9231
$ws_CUL.Range("L122").Value = 9980.3079
$ws_CUL.Range("N122").Value = -14880.3079

# GSM row 56
$ws_GSM.Range("H56").Value = 6000
$ws_GSM.Range("I56").Value = 6000
$ws_GSM.Range("J56").Value = 0
$ws_GSM.Range("K56").Value = 6000
$ws_GSM.Range("L56").Value = 0
$ws_GSM.Range("M56").Value = -5248
$ws_GSM.Range("N56").ClearContents()

# GSM row 80
$ws_GSM.Range("H80").Value = 4668.8335
$ws_GSM.Range("I80").Value = 4760.857
$ws_GSM.Range("J80").Value = 4540
$ws_GSM.Range("K80").Value = 4760.857
$ws_GSM.Range("L80").Value = 4540
$ws_GSM.Range("M80").Value = -3762.857
$ws_GSM.Range("N80").Value = -6536

# GSM row 83
$ws_GSM.Range("H83").Value = 4668.8335
$ws_GSM.Range("I83").Value = 4760.857
$ws_GSM.Range("J83").Value = 4540
$ws_GSM.Range("K83").Value = 23804.285
$ws_GSM.Range("L83").Value = 22700
$ws_GSM.Range("M83").Value = -18812.285
$ws_GSM.Range("N83").Value = -32684

# GSM row 132
$ws_GSM.Range("H132").Value = 2480.162
$ws_GSM.Range("I132").Value = 2241.2322
$ws_GSM.Range("K132").Value = 6723.696599999999
$ws_GSM.Range("M132").Value = -4193.696599999999

# LTW row 7
$ws_LTW.Range("H7").Value = 3852.3784
$ws_LTW.Range("I7").Value = 2743.5173
$ws_LTW.Range("J7").Value = 7872
$ws_LTW.Range("K7").Value = 2743.5173
$ws_LTW.Range("L7").Value = 7872
$ws_LTW.Range("M7").Value = -2631.5173
$ws_LTW.Range("N7").Value = -8096

# LTW row 53
$ws_LTW.Range("H53").Value = 15499
$ws_LTW.Range("I53").Value = 14999
$ws_LTW.Range("J53").Value = 15999
$ws_LTW.Range("K53").Value = 14999
$ws_LTW.Range("L53").Value = 15999
$ws_LTW.Range("M53").Value = -14481
$ws_LTW.Range("N53").Value = -17035

# LTW row 100
$ws_LTW.Range("H100").Value = 3928.5715
$ws_LTW.Range("I100").Value = 3500
$ws_LTW.Range("K100").Value = 3500
$ws_LTW.Range("M100").Value = -2959

# LTW row 122
$ws_LTW.Range("H122").Value = 5600.5713
$ws_LTW.Range("I122").Value = 1835
$ws_LTW.Range("K122").Value = 5505
$ws_LTW.Range("M122").Value = -3055

# LTW row 126
$ws_LTW.Range("H126").Value = 3852.3784
$ws_LTW.Range("I126").Value = 2743.5173
$ws_LTW.Range("J126").Value = 7872
$ws_LTW.Range("K126").Value = 8230.5519
$ws_LTW.Range("L126").Value = 23616
$ws_LTW.Range("M126").Value = -5760.5519
$ws_LTW.Range("N126").Value = -28556

# LTW row 132
$ws_LTW.Range("H132").Value = 3114.276
$ws_LTW.Range("I132").Value = 3110.111
$ws_LTW.Range("J132").Value = 3128.6924
$ws_LTW.Range("K132").Value = 9330.332999999999
$ws_LTW.Range("L132").Value = 9386.0772
$ws_LTW.Range("M132").Value = -6800.332999999999
$ws_LTW.Range("N132").Value = -14446.0772

# LTW row 136
$ws_LTW.Range("H136").Value = 4628.795
$ws_LTW.Range("I136").Value = 4228.357
$ws_LTW.Range("J136").Value = 5648.091
$ws_LTW.Range("K136").Value = 12685.071
$ws_LTW.Range("L136").Value = 16944.273
$ws_LTW.Range("M136").Value = -10135.071
$ws_LTW.Range("N136").Value = -22044.273

# WVR row 113
$ws_WVR.Range("H113").Value = 3626453.8
$ws_WVR.Range("I113").Value = 3971802.5
$ws_WVR.Range("J113").Value = 292.5
$ws_WVR.Range("K113").Value = 11915407.5
$ws_WVR.Range("L113").Value = 877.5
$ws_WVR.Range("M113").Value = -11913237.5
$ws_WVR.Range("N113").Value = -5217.5

# WVR row 122
$ws_WVR.Range("H122").Value = 2712.353
$ws_WVR.Range("I122").Value = 2718.5386
$ws_WVR.Range("K122").Value = 8155.6158
$ws_WVR.Range("M122").Value = -5705.6158

# WVR row 126
$ws_WVR.Range("H126").Value = 2410.35
$ws_WVR.Range("I126").Value = 2122.611
$ws_WVR.Range("K126").Value = 6367.833
$ws_WVR.Range("M126").Value = -3897.833

# WVR row 132
$ws_WVR.Range("H132").Value = 2067.9565
$ws_WVR.Range("I132").Value = 2008.2
$ws_WVR.Range("K132").Value = 6024.6
$ws_WVR.Range("M132").Value = -3494.6

# WVR row 136
$ws_WVR.Range("H136").Value = 14655.125
$ws_WVR.Range("I136").Value = 14655.125
$ws_WVR.Range("J136").Value = 0
$ws_WVR.Range("K136").Value = 43965.375
$ws_WVR.Range("L136").Value = 0
$ws_WVR.Range("M136").Value = -41415.375
$ws_WVR.Range("N136").ClearContents()
